# v.0.5.0 One more instruction working
#
# - Rename the sheet "Instruction format" -> "Instructions format"
# - Add a new title row above the existing table: B2 = "JACA-2 Instructions
#   Format", bold (new font-only cell style), no fill/border.
# - Move the active selection from B4:G4 to the new title cell B3's row
#   (single cell B3, matching the author's post-edit selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet.
$ws.Name = "Instructions format"

# New title cell above the existing content.
$ws.Range("B2").Value = "JACA-2 Instructions Format"
$ws.Range("B2").Font.Bold = $true

# Update the selection to match the saved workbook state.
$ws.Range("B3").Select()
